$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Paragraph "We strive to allow the community to control the fate of
# the coin..." -> "The self funding mechanism that pays for
# development and encourages the community to decide how budget funds
# are spent." (SmartHive section, new formatting: Times New Roman /
# Open Sans, smaller sizes, explicit shading on the run).
# ---------------------------------------------------------------------
$r1 = $d.Content
$null = $r1.Find.Execute("We strive to", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target1 = $r1.Paragraphs(1).Range
$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/><w:b/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">The self funding mechanism that pays for development and encourages the community to decide how budget funds are spent.</w:t></w:r></w:p>
'@
$target1.InsertXML($xml1)

# ---------------------------------------------------------------------
# Paragraph "SmartCash mining prevents mining centralization..." ->
# "SmartCash now has SmartMining that prevents mining attacks..."
# (THE SMARTCASH MINING CONCEPT section; paragraph mark reformatted to
# Times New Roman 12pt, shading/border dropped from the pPr rPr).
# ---------------------------------------------------------------------
$r2 = $d.Content
$null = $r2.Find.Execute("SmartCash mining prevents", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target2 = $r2.Paragraphs(1).Range
$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t xml:space="preserve">SmartCash now has SmartMining that prevents mining attacks. Mining can be done by anyone with a computer with one or more graphics cards. ASICs have yet to be created for the Keccak mining algorithm and it’s probably safe to assume no ASICs will be created for quite some time.</w:t></w:r></w:p>
'@
$target2.InsertXML($xml2)

Write-Host "Applied Welcome.docx (German) translation updates."
